# "basic risk analysis added"
# Update the risk-analysis data point and refresh the active selection,
# matching the cell that was last reviewed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Ark1")

# Risk value bump: D5 (goo / dk) from 85 -> 500
$ws.Range("D5").Value = 500

# Leave the selection where the analyst left off after editing
$ws.Range("F7").Select()
